# "icons on each button vers2"
#
# Adds a new worksheet named "icons" as the last sheet in the workbook,
# puts the label "hdd.png" in A2/B2 (this introduces a new shared string),
# formats the small "staircase" block of cells A2:D3 / B4,D4 / B5,D5 / D6
# (mirrors the used-range footprint the original author ended up with),
# and makes "icons" the active/selected sheet (so workbook.xml's
# activeTab moves from the "color" sheet to this new last sheet, and the
# previously tabSelected sheet reverts to unselected).

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the very end of the tab strip -------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "icons"

# --- content ----------------------------------------------------------------
$ws.Range("A2").Value = "hdd.png"
$ws.Range("B2").Value = "hdd.png"

# --- formatting footprint (reuses the existing "General" number format,
#     so no new style entries are introduced) -------------------------------
$ws.Range("A2:D3").NumberFormat = "General"
$ws.Range("B4").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("B5").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "General"

# --- make "icons" the active sheet / cell selection -------------------------
$ws.Activate()
[void]$ws.Range("D16").Select()
